$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header update ---
# Existing C1 ("WC48 P5F") shifts right to D1; new headers are inserted
# into C1, E1, F1, G1.
$ws.Range("C1").Value = "WC47 NACP"
$ws.Range("D1").Value = "WC48 P5F"
$ws.Range("E1").Value = "WC49 P5H"
$ws.Range("F1").Value = "WV50 FILTER"
$ws.Range("G1").Value = "SPL"

# --- New ticket rows 13-26 ---
$data = @(
    @("2024-05-15","11:14:51","-","Cámara no detecta Pcb","-","-","-"),
    @("2024-05-15","11:14:58","-","Cámara no detecta Top cover","-","-","-"),
    @("2024-05-15","11:15:32","-","Tornillo atascado en tolva","-","-","-"),
    @("2024-05-15","11:15:35","-","No coloca bien el sealling","-","-","-"),
    @("2024-05-15","11:16:16","-","-","-","-","Colisión placas"),
    @("2024-05-15","11:16:21","-","-","-","-","Error en sensor de salida"),
    @("2024-05-15","11:16:23","-","-","-","-","Marco atascado en parte inferior"),
    @("2024-05-15","11:16:32","-","-","-","-","Colisión placas"),
    @("2024-05-15","11:24:41","-","-","-","Cover atascado","-"),
    @("2024-05-15","11:24:44","-","-","-","No coloca bien la pcb","-"),
    @("2024-05-15","11:24:48","-","-","-","NOK Soldadura Plástico","-"),
    @("2024-05-15","11:24:50","-","-","-","Fallo cámara cover","-"),
    @("2024-05-15","11:24:53","-","-","-","Fallo cámara QR","-"),
    @("2024-05-15","11:24:55","-","-","-","No coloca bien foam","-")
)

$startRow = 13
$endRow = $startRow + $data.Count - 1

# Column A holds date-looking strings ("2024-05-15"). Force the whole
# column range to be stored as literal text (matching the original rows'
# representation) instead of letting Excel auto-convert it to a date
# serial number; the style is reset to Normal afterwards so no visible
# formatting change is left behind.
$ws.Range("A$startRow" + ":A$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

$ws.Range("A$startRow" + ":A$endRow").Style = "Normal"
